$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.0430770272853
$ws.Range("C2").Value = 9.33960046570504
$ws.Range("E2").Value = 11.33263631685594
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.623531093483944
$ws.Range("I2").Value = 21.02444577804193
$ws.Range("L2").Value = 9.779408433516078
$ws.Range("N2").Value = 16.9816204737432
$ws.Range("O2").Value = 20.96264503873406
$ws.Range("B3").Value = 14.53230061879652
$ws.Range("C3").Value = 9.131815313324816
$ws.Range("E3").Value = 11.36676225565127
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.625505936589418
$ws.Range("I3").Value = 21.13184774284021
$ws.Range("L3").Value = 9.753841679127893
$ws.Range("N3").Value = 17.02230160730527
$ws.Range("O3").Value = 21.0265966940832
$ws.Range("B4").Value = 14.211065880304
$ws.Range("C4").Value = 9.000689545531557
$ws.Range("E4").Value = 11.38962252686372
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.62678313366008
$ws.Range("I4").Value = 21.20264681174368
$ws.Range("L4").Value = 9.739849734517943
$ws.Range("N4").Value = 17.04911196219501
$ws.Range("O4").Value = 21.07155023170489
$ws.Range("B5").Value = 14.07844346458596
$ws.Range("C5").Value = 8.946410638612825
$ws.Range("E5").Value = 11.39941759420968
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.627319906050325
$ws.Range("I5").Value = 21.23271646862696
$ws.Range("L5").Value = 9.734581095834665
$ws.Range("N5").Value = 17.06049879229157
$ws.Range("O5").Value = 21.09129398156967
$ws.Range("B6").Value = 14.0563240782176
$ws.Range("C6").Value = 8.937348121460062
$ws.Range("E6").Value = 11.40107300257407
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.627410022934688
$ws.Range("I6").Value = 21.23778304467472
$ws.Range("L6").Value = 9.733732520409013
$ws.Range("N6").Value = 17.06241745535505
$ws.Range("O6").Value = 21.09465831362905
$ws.Range("B7").Value = 14.20928396347803
$ws.Range("C7").Value = 8.99996087443156
$ws.Range("E7").Value = 11.38975268602951
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.626790306677337
$ws.Range("I7").Value = 21.20304741086784
$ws.Range("L7").Value = 9.739776920612304
$ws.Range("N7").Value = 17.04926365983754
$ws.Range("O7").Value = 21.07181074091164
$ws.Range("B8").Value = 14.86865559534553
$ws.Range("C8").Value = 9.268719394720241
$ws.Range("E8").Value = 11.34400709055803
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.624198635116394
$ws.Range("I8").Value = 21.0604693084837
$ws.Range("L8").Value = 9.770241518882255
$ws.Range("N8").Value = 16.99526753728933
$ws.Range("O8").Value = 20.98351219075155
$ws.Range("B9").Value = 16.0934298020004
$ws.Range("C9").Value = 9.765692703681941
$ws.Range("E9").Value = 11.26943662352095
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.619626894953373
$ws.Range("I9").Value = 20.81948061616312
$ws.Range("L9").Value = 9.843320833995568
$ws.Range("N9").Value = 16.90388516647577
$ws.Range("O9").Value = 20.85570349205033
$ws.Range("B10").Value = 16.94240995386489
$ws.Range("C10").Value = 10.11010613831209
$ws.Range("E10").Value = 11.22388614205287
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.616575997352896
$ws.Range("I10").Value = 20.6660913800636
$ws.Range("L10").Value = 9.904854000866067
$ws.Range("N10").Value = 16.84554532185323
$ws.Range("O10").Value = 20.78972918766196
$ws.Range("B11").Value = 17.31596444255711
$ws.Range("C11").Value = 10.26184373675861
$ws.Range("E11").Value = 11.20517129112285
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.61525424077359
$ws.Range("I11").Value = 20.60147830740542
$ws.Range("L11").Value = 9.934480713797431
$ws.Range("N11").Value = 16.82090670673785
$ws.Range("O11").Value = 20.76583338523513
$ws.Range("B12").Value = 17.45548572853392
$ws.Range("C12").Value = 10.31855913216458
$ws.Range("E12").Value = 11.19837309387291
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.614763180146083
$ws.Range("I12").Value = 20.57775626697595
$ws.Range("L12").Value = 9.945928295930614
$ws.Range("N12").Value = 16.81184932922546
$ws.Range("O12").Value = 20.75766794490472
$ws.Range("B13").Value = 17.42552512074742
$ws.Range("C13").Value = 10.30637807811003
$ws.Range("E13").Value = 11.19982436370518
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.61486851885013
$ws.Range("I13").Value = 20.58283202873206
$ws.Range("L13").Value = 9.943452794844234
$ws.Range("N13").Value = 16.81378787640839
$ws.Range("O13").Value = 20.75938717269594
$ws.Range("B14").Value = 17.32748235460506
$ws.Range("C14").Value = 10.26652485066248
$ws.Range("E14").Value = 11.20460621203114
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.615213651597124
$ws.Range("I14").Value = 20.5995117226489
$ws.Range("L14").Value = 9.935417970451427
$ws.Range("N14").Value = 16.82015608802272
$ws.Range("O14").Value = 20.76514388455486
$ws.Range("B15").Value = 17.26717302637661
$ws.Range("C15").Value = 10.24201569299484
$ws.Range("E15").Value = 11.20757283668849
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.615426285942971
$ws.Range("I15").Value = 20.60982569310928
$ws.Range("L15").Value = 9.93052597874123
$ws.Range("N15").Value = 16.82409230296897
$ws.Range("O15").Value = 20.76878518701408
$ws.Range("B16").Value = 16.9177324389945
$ws.Range("C16").Value = 10.10008737309597
$ws.Range("E16").Value = 11.22514959156062
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.61666370357758
$ws.Range("I16").Value = 20.67041816376677
$ws.Range("L16").Value = 9.902950175356604
$ws.Range("N16").Value = 16.84719370795365
$ws.Range("O16").Value = 20.79141424946715
$ws.Range("B17").Value = 16.70003434178127
$ws.Range("C17").Value = 10.01172896806899
$ws.Range("E17").Value = 11.23644637534291
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.617439718586707
$ws.Range("I17").Value = 20.70891445865095
$ws.Range("L17").Value = 9.886447532080826
$ws.Range("N17").Value = 16.86185202938864
$ws.Range("O17").Value = 20.80686584837446
$ws.Range("B18").Value = 16.5736354232893
$ws.Range("C18").Value = 9.960444651035365
$ws.Range("E18").Value = 11.24313281703842
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.617892287454906
$ws.Range("I18").Value = 20.73154243654056
$ws.Range("L18").Value = 9.877109949818212
$ws.Range("N18").Value = 16.87046202188517
$ws.Range("O18").Value = 20.81632868262179
$ws.Range("B19").Value = 16.53063922576298
$ws.Range("C19").Value = 9.943002249785081
$ws.Range("E19").Value = 11.2454291564409
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.618046590180098
$ws.Range("I19").Value = 20.73928725258417
$ws.Range("L19").Value = 9.873975099712407
$ws.Range("N19").Value = 16.87340796314074
$ws.Range("O19").Value = 20.81963135077221
$ws.Range("B20").Value = 16.7233321797827
$ws.Range("C20").Value = 10.02118303673868
$ws.Range("E20").Value = 11.23522426861227
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.617356466483491
$ws.Range("I20").Value = 20.70476614995869
$ws.Range("L20").Value = 9.888188343190825
$ws.Range("N20").Value = 16.8602731132042
$ws.Range("O20").Value = 20.80516140954734
$ws.Range("B21").Value = 17.35633328503975
$ws.Range("C21").Value = 10.27825117025611
$ws.Range("E21").Value = 11.20319383018117
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.615112021375983
$ws.Range("I21").Value = 20.59459223553901
$ws.Range("L21").Value = 9.937771842799826
$ws.Range("N21").Value = 16.81827819240242
$ws.Range("O21").Value = 20.76342899520638
$ws.Range("B22").Value = 17.75870898664012
$ws.Range("C22").Value = 10.44190751399667
$ws.Range("E22").Value = 11.18394299562093
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.613700264793505
$ws.Range("I22").Value = 20.52693409898021
$ws.Range("L22").Value = 9.971506677388605
$ws.Range("N22").Value = 16.79242152348837
$ws.Range("O22").Value = 20.74130497319281
$ws.Range("B23").Value = 17.54502486749947
$ws.Range("C23").Value = 10.35496991387575
$ws.Range("E23").Value = 11.19406347854539
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.614448717954813
$ws.Range("I23").Value = 20.56264581181028
$ws.Range("L23").Value = 9.953382376946697
$ws.Range("N23").Value = 16.80607646001644
$ws.Range("O23").Value = 20.75264053533171
$ws.Range("B24").Value = 16.71280308483872
$ws.Range("C24").Value = 10.0169103631509
$ws.Range("E24").Value = 11.23577618574278
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.617394084734266
$ws.Range("I24").Value = 20.70664005567392
$ws.Range("L24").Value = 9.887400854744401
$ws.Range("N24").Value = 16.86098637197196
$ws.Range("O24").Value = 20.80593018149931
$ws.Range("B25").Value = 15.77042376950782
$ws.Range("C25").Value = 9.634724144620135
$ws.Range("E25").Value = 11.2879884706404
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.620809356195789
$ws.Range("I25").Value = 20.88052818776773
$ws.Range("L25").Value = 9.822151421787503
$ws.Range("N25").Value = 16.92705829865755
$ws.Range("O25").Value = 20.85570349205033
